$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 432 - this shifts the existing rows
# 432..493 down to 433..494 (matching the diff's dimension change
# from A1:R493 to A1:R494 and the observed row-by-row downward shift).
$ws.Rows.Item(432).Insert()

# Populate the newly inserted row 432 with a new weekly record for
# "Feria Lagunitas de Puerto Montt" / Cebollín, identical to the
# (now shifted) row below except for the date (D) and volume (J).
$ws.Range("A432").Value = 4
$ws.Range("B432").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C432").Value = "Los Lagos"
$ws.Range("D432").Value = 45142
$ws.Range("E432").Value = 10
$ws.Range("F432").Value = 100112037
$ws.Range("G432").Value = "Cebollín"
$ws.Range("H432").Value = "Sin especificar"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 180
$ws.Range("K432").Value = 6000
$ws.Range("L432").Value = 6000
$ws.Range("M432").Value = 6000
$ws.Range("N432").Value = "$/paquete 36 unidades"
$ws.Range("O432").Value = "Región Metropolitana"
$ws.Range("P432").Value = 167
$ws.Range("Q432").Value = 36
$ws.Range("R432").Value = "Hortaliza"
